$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MCH230-1 entry
$ws.Range("A2").Value = "MCH230-1"
$ws.Range("C2").Value = "VUKANI BASEBENZI- WORKERS WORLD"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"

# Row 3: MCH230-2 entry
$ws.Range("A3").Value = "MCH230-2"
$ws.Range("C3").Value = "GERMANY ANTI-APARTHEID GROUPS PUBLICATION"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24G | GRAP COUNT NUMER: NONE"
